$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2543679979581253
$ws.Range("D2").Value = 0.1203927657961934
$ws.Range("E2").Value = 0.1417083679642417
$ws.Range("F2").Value = 1.225464199728314
$ws.Range("G2").Value = 0.7347609656863057
$ws.Range("H2").Value = 0.8372351477910271
$ws.Range("J2").Value = 0.2181232318074393
$ws.Range("L2").Value = 0.2282406447922511
$ws.Range("M2").Value = 11.06113479499805
$ws.Range("O2").Value = 3.149964013821844
$ws.Range("C3").Value = 0.264913751971017
$ws.Range("D3").Value = 0.1232911994105805
$ws.Range("E3").Value = 0.1411351465171471
$ws.Range("F3").Value = 1.263195159704452
$ws.Range("G3").Value = 0.7561250787211264
$ws.Range("H3").Value = 0.8566154120698855
$ws.Range("J3").Value = 0.2111458522780225
$ws.Range("L3").Value = 0.2138939518756047
$ws.Range("M3").Value = 9.740127417205997
$ws.Range("O3").Value = 3.235294182199283
$ws.Range("C4").Value = 0.2717416006280207
$ws.Range("D4").Value = 0.1251729111397282
$ws.Range("E4").Value = 0.1409471580577275
$ws.Range("F4").Value = 1.287933956415817
$ws.Range("G4").Value = 0.7703908290537953
$ws.Range("H4").Value = 0.8693078194938977
$ws.Range("J4").Value = 0.207109829528747
$ws.Range("L4").Value = 0.2052145579937559
$ws.Range("M4").Value = 8.925350066039641
$ws.Range("O4").Value = 3.291765736595082
$ws.Range("C5").Value = 0.2746116862619594
$ws.Range("D5").Value = 0.1259652117347017
$ws.Range("E5").Value = 0.1409111362799429
$ws.Range("F5").Value = 1.298406112951742
$ws.Range("G5").Value = 0.7764890881850306
$ws.Range("H5").Value = 0.8746778109407529
$ws.Range("J5").Value = 0.2055261824246202
$ws.Range("L5").Value = 0.2017095315205637
$ws.Range("M5").Value = 8.592380070635841
$ws.Range("O5").Value = 3.315792328691956
$ws.Range("C6").Value = 0.2750935281950913
$ws.Range("D6").Value = 0.1260983069565782
$ws.Range("E6").Value = 0.1409075878195729
$ws.Range("F6").Value = 1.300168488673581
$ws.Range("G6").Value = 0.7775187914958153
$ws.Range("H6").Value = 0.8755813886836989
$ws.Range("J6").Value = 0.2052668698798925
$ws.Range("L6").Value = 0.2011294327255655
$ws.Range("M6").Value = 8.537033340657558
$ws.Range("O6").Value = 3.319842829363139
$ws.Range("C7").Value = 0.2717799539108192
$ws.Range("D7").Value = 0.1251834934136689
$ws.Range("E7").Value = 0.1409465087590434
$ws.Range("F7").Value = 1.28807361058864
$ws.Range("G7").Value = 0.7704719240267082
$ws.Range("H7").Value = 0.8693794426196106
$ws.Range("J7").Value = 0.2070882263566745
$ws.Range("L7").Value = 0.2051671596696991
$ws.Range("M7").Value = 8.920863349580259
$ws.Range("O7").Value = 3.292085677330007
$ws.Range("C8").Value = 0.2579303589089825
$ws.Range("D8").Value = 0.1213708558229918
$ws.Range("E8").Value = 0.141476288942961
$ws.Range("F8").Value = 1.238144980829873
$ws.Range("G8").Value = 0.7418867053430134
$ws.Range("H8").Value = 0.8437519255625645
$ws.Range("J8").Value = 0.2156651107155341
$ws.Range("L8").Value = 0.2232665864986387
$ws.Range("M8").Value = 10.60640150939537
$ws.Range("O8").Value = 3.17853238830665
$ws.Range("C9").Value = 0.2336151136991784
$ws.Range("D9").Value = 0.1147120679169333
$ws.Range("E9").Value = 0.1438472777948689
$ws.Range("F9").Value = 1.152900959937419
$ws.Range("G9").Value = 0.695113425370991
$ws.Range("H9").Value = 0.7998598359834688
$ws.Range("J9").Value = 0.2345177408085561
$ws.Range("L9").Value = 0.2598235240090361
$ws.Range("M9").Value = 13.88375980679825
$ws.Range("O9").Value = 2.988733909757343
$ws.Range("C10").Value = 0.2175463572684411
$ws.Range("D10").Value = 0.110329477694826
$ws.Range("E10").Value = 0.14644711647977
$ws.Range("F10").Value = 1.098259465981968
$ws.Range("G10").Value = 0.6666457784237991
$ws.Range("H10").Value = 0.7715904057251493
$ws.Range("J10").Value = 0.2497052256671708
$ws.Range("L10").Value = 0.2873891985618968
$ws.Range("M10").Value = 16.27673002307228
$ws.Range("O10").Value = 2.870045084558967
$ws.Range("C11").Value = 0.2106398318066667
$ws.Range("D11").Value = 0.1084488747237629
$ws.Range("E11").Value = 0.1478268348604317
$ws.Range("F11").Value = 1.07519396195255
$ws.Range("G11").Value = 0.6550275230294176
$ws.Range("H11").Value = 0.7596151657865988
$ws.Range("J11").Value = 0.2569274584305958
$ws.Range("L11").Value = 0.3000976987547119
$ws.Range("M11").Value = 17.36268695714216
$ws.Range("O11").Value = 2.820712792546033
$ws.Range("C12").Value = 0.2080836756379103
$ws.Range("D12").Value = 0.10775321770506
$ws.Range("E12").Value = 0.1483785626371628
$ws.Range("F12").Value = 1.066722157217015
$ws.Range("G12").Value = 0.6508238966680864
$ws.Range("H12").Value = 0.7552094968481811
$ws.Range("J12").Value = 0.2597093706285989
$ws.Range("L12").Value = 0.3049356005263348
$ws.Range("M12").Value = 17.77358224885535
$ws.Range("O12").Value = 2.80271516835839
$ws.Range("C13").Value = 0.2086315373156165
$ws.Range("D13").Value = 0.1079023029251687
$ws.Range("E13").Value = 0.1482584216911036
$ws.Range("F13").Value = 1.068534951733376
$ws.Range("G13").Value = 0.651720436183794
$ws.Range("H13").Value = 0.7561525647435445
$ws.Range("J13").Value = 0.2591081131475619
$ws.Range("L13").Value = 0.3038925188847088
$ws.Range("M13").Value = 17.68510266560077
$ws.Range("O13").Value = 2.806560667355711
$ws.Range("C14").Value = 0.2104283419113386
$ws.Range("D14").Value = 0.108391310851772
$ws.Range("E14").Value = 0.1478716339943702
$ws.Range("F14").Value = 1.074491689470776
$ws.Range("G14").Value = 0.6546777352422311
$ws.Range("H14").Value = 0.7592501106450698
$ws.Range("J14").Value = 0.2571553743478034
$ws.Range("L14").Value = 0.3004951981750992
$ws.Range("M14").Value = 17.3964979688871
$ws.Range("O14").Value = 2.819218342665692
$ws.Range("C15").Value = 0.2115366819537723
$ws.Range("D15").Value = 0.108692996273291
$ws.Range("E15").Value = 0.1476385539029863
$ws.Range("F15").Value = 1.078174712396255
$ws.Range("G15").Value = 0.6565148184352836
$ws.Range("H15").Value = 0.7611643151192311
$ws.Range("J15").Value = 0.2559654470877746
$ws.Range("L15").Value = 0.2984175952995827
$ws.Range("M15").Value = 17.21967725761499
$ws.Range("O15").Value = 2.827060939022516
$ws.Range("C16").Value = 0.2180059319145133
$ws.Range("D16").Value = 0.1104546739433943
$ws.Range("E16").Value = 0.146360991517696
$ws.Range("F16").Value = 1.099803322401115
$ws.Range("G16").Value = 0.6674322363309955
$ws.Range("H16").Value = 0.7723909777404572
$ws.Range("J16").Value = 0.2492397055779492
$ws.Range("L16").Value = 0.2865621643597365
$ws.Range("M16").Value = 16.20571111070853
$ws.Range("O16").Value = 2.873363971608768
$ws.Range("C17").Value = 0.2220787378908593
$ws.Range("D17").Value = 0.1115645236992222
$ws.Range("E17").Value = 0.1456283577197226
$ws.Range("F17").Value = 1.113534094591898
$ws.Range("G17").Value = 0.6744739668517923
$ws.Range("H17").Value = 0.7795060859307483
$ws.Range("J17").Value = 0.2451952941782025
$ws.Range("L17").Value = 0.2793332855356994
$ws.Range("M17").Value = 15.58303467738915
$ws.Range("O17").Value = 2.902972325673801
$ws.Range("C18").Value = 0.2244591883632641
$ws.Range("D18").Value = 0.1122135074801776
$ws.Range("E18").Value = 0.1452254542252334
$ws.Range("F18").Value = 1.121599961183563
$ws.Range("G18").Value = 0.6786492715534962
$ws.Range("H18").Value = 0.7836816630338177
$ws.Range("J18").Value = 0.2428984285612046
$ws.Range("L18").Value = 0.2751912299590202
$ws.Range("M18").Value = 15.22463842385133
$ws.Range("O18").Value = 2.920440064187545
$ws.Range("C19").Value = 0.2252716374733446
$ws.Range("D19").Value = 0.1124350606765958
$ws.Range("E19").Value = 0.1450921859565923
$ws.Range("F19").Value = 1.124359676848506
$ws.Range("G19").Value = 0.680084307708853
$ws.Range("H19").Value = 0.7851096667319268
$ws.Range("J19").Value = 0.2421257365457024
$ws.Range("L19").Value = 0.2737914803823571
$ws.Range("M19").Value = 15.10324756344653
$ws.Range("O19").Value = 2.926429128266648
$ws.Range("C20").Value = 0.2216412516024544
$ws.Range("D20").Value = 0.1114452767547363
$ws.Range("E20").Value = 0.1457044278448762
$ws.Range("F20").Value = 1.112054976935561
$ws.Range("G20").Value = 0.6737113845968423
$ws.Range("H20").Value = 0.778740050746535
$ws.Range("J20").Value = 0.2456227749350717
$ws.Range("L20").Value = 0.2801011679217424
$ws.Range("M20").Value = 15.64934514350733
$ws.Range("O20").Value = 2.899775059452367
$ws.Range("C21").Value = 0.2098989605222439
$ws.Range("D21").Value = 0.1082472280690823
$ws.Range("E21").Value = 0.1479844412499602
$ws.Range("F21").Value = 1.072734883788996
$ws.Range("G21").Value = 0.6538037509078549
$ws.Range("H21").Value = 0.7583367676183457
$ws.Range("J21").Value = 0.2577276491907128
$ws.Range("L21").Value = 0.301492372328994
$ws.Range("M21").Value = 17.48127683314266
$ws.Range("O21").Value = 2.815481820563832
$ws.Range("C22").Value = 0.2025703106068164
$ws.Range("D22").Value = 0.1062533260714531
$ws.Range("E22").Value = 0.1496455649023005
$ws.Range("F22").Value = 1.048569963187688
$ws.Range("G22").Value = 0.6419372503052472
$ws.Range("H22").Value = 0.7457553963112815
$ws.Range("J22").Value = 0.265913872165001
$ws.Range("L22").Value = 0.3156218395873225
$ws.Range("M22").Value = 18.67663252436517
$ws.Range("O22").Value = 2.764381307020443
$ws.Range("C23").Value = 0.2064497179124736
$ws.Range("D23").Value = 0.1073086284999931
$ws.Range("E23").Value = 0.1487430273722978
$ws.Range("F23").Value = 1.061325280012191
$ws.Range("G23").Value = 0.6481644359218848
$ws.Range("H23").Value = 0.7524007507230124
$ws.Range("J23").Value = 0.2615188979006291
$ws.Range("L23").Value = 0.3080666130599639
$ws.Range("M23").Value = 18.038808568028
$ws.Range("O23").Value = 2.791285057404991
$ws.Range("C24").Value = 0.2218389179560081
$ws.Range("D24").Value = 0.1114991543135133
$ws.Range("E24").Value = 0.1456699796182974
$ws.Range("F24").Value = 1.11272315109391
$ws.Range("G24").Value = 0.6740557530419977
$ws.Range("H24").Value = 0.7790861106191898
$ws.Range("J24").Value = 0.2454294228194982
$ws.Range("L24").Value = 0.2797539650174485
$ws.Range("M24").Value = 15.61936746602532
$ws.Range("O24").Value = 2.901219156883201
$ws.Range("C25").Value = 0.2398822515964127
$ws.Range("D25").Value = 0.1164247746180322
$ws.Range("E25").Value = 0.1430586731072836
$ws.Range("F25").Value = 1.174578481603895
$ws.Range("G25").Value = 0.7067487000100101
$ws.Range("H25").Value = 0.8110424122885433
$ws.Range("J25").Value = 0.2291902623309028
$ws.Range("L25").Value = 0.2498142609829159
$ws.Range("M25").Value = 12.99993745599107
$ws.Range("O25").Value = 3.036486409104583
